# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect the repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -1
    5  = 5
    7  = 8
    10 = 5
    11 = 5
    12 = -2
    13 = -2
    17 = -7
    18 = 0
    19 = 4
    21 = -6
    22 = -2
    23 = -4
    25 = 2
    26 = -5
    27 = -1
    29 = -12
    31 = -2
    32 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
